$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2..189).
# The commit updates that date from 2023-09-03 (45172) to 2023-09-06 (45175)
# for all rows, leaving everything else untouched.
$newDate = 45175
$lastRow = 189

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
